# Totals/ml parameter search script: swap the roster rows for
# Jaylen Brown (row 7) and Malcolm Brogdon (row 8) so each player's stat
# columns (No., Pos, Ht, Wt, Birth Date, College, bbref url) line up with
# the other player. Columns that already hold identical values in both
# rows (Exp, country) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 7
$row2 = 8

# Columns that actually differ between the two rows and need to swap.
$cols = @(2, 3, 4, 5, 6, 7, 10, 11)   # B, C, D, E, F, G, J, K

foreach ($c in $cols) {
    $cell1 = $ws.Cells.Item($row1, $c)
    $cell2 = $ws.Cells.Item($row2, $c)

    $v1 = $cell1.Value()
    $v2 = $cell2.Value()

    $cell1.Value = $v2
    $cell2.Value = $v1
}
